$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 201.4397426666667
$ws.Range("H2").Value = 604.3192280000001
$ws.Range("I2").Value = 0.4833500233086392
$ws.Range("J2").Value = 0.4833500233086393
$ws.Range("M2").Value = 0.74396
$ws.Range("N2").Value = 2.23188
$ws.Range("O2").Value = 0.006259003216804254
$ws.Range("P2").Value = 0.006259003216804255
$ws.Range("Q2").Value = 149.8631109542933
$ws.Range("R2").Value = 1348.76799858864
$ws.Range("S2").Value = 0.003025289350731184
$ws.Range("T2").Value = 0.003025289350731185
$ws.Range("G3").Value = 201.4397426666667
$ws.Range("H3").Value = 604.3192280000001
$ws.Range("I3").Value = 0.4833500233086392
$ws.Range("J3").Value = 0.4833500233086393
$ws.Range("M3").Value = 88.14978533333333
$ws.Range("O3").Value = 0.7416121699579786
$ws.Range("P3").Value = 0.7416121699579786
$ws.Range("Q3").Value = 17756.87007366858
$ws.Range("R3").Value = 159811.8306630172
$ws.Range("S3").Value = 0.3584582596351595
$ws.Range("T3").Value = 0.3584582596351595
$ws.Range("G4").Value = 201.4397426666667
$ws.Range("H4").Value = 604.3192280000001
$ws.Range("I4").Value = 0.4833500233086392
$ws.Range("J4").Value = 0.4833500233086393
$ws.Range("M4").Value = 29.76859933333333
$ws.Range("N4").Value = 89.305798
$ws.Range("O4").Value = 0.2504459365921425
$ws.Range("P4").Value = 0.2504459365921425
$ws.Range("Q4").Value = 5996.578989253771
$ws.Range("R4").Value = 53969.21090328395
$ws.Range("S4").Value = 0.1210530492893661
$ws.Range("T4").Value = 0.1210530492893661
$ws.Range("G5").Value = 201.4397426666667
$ws.Range("H5").Value = 604.3192280000001
$ws.Range("I5").Value = 0.4833500233086392
$ws.Range("J5").Value = 0.4833500233086393
$ws.Range("M5").Value = 0.2000323333333334
$ws.Range("N5").Value = 0.6000970000000001
$ws.Range("O5").Value = 0.00168289023307462
$ws.Range("P5").Value = 0.00168289023307462
$ws.Range("Q5").Value = 40.29446175167957
$ws.Range("R5").Value = 362.6501557651161
$ws.Range("S5").Value = 0.000813425033382499
$ws.Range("T5").Value = 0.0008134250333824992
$ws.Range("I6").Value = 0.1569674599353791
$ws.Range("J6").Value = 0.1569674599353792
$ws.Range("M6").Value = 0.74396
$ws.Range("N6").Value = 2.23188
$ws.Range("O6").Value = 0.006259003216804254
$ws.Range("P6").Value = 0.006259003216804255
$ws.Range("Q6").Value = 48.66790261741333
$ws.Range("R6").Value = 438.01112355672
$ws.Range("S6").Value = 0.0009824598366691309
$ws.Range("T6").Value = 0.0009824598366691313
$ws.Range("I7").Value = 0.1569674599353791
$ws.Range("J7").Value = 0.1569674599353792
$ws.Range("M7").Value = 88.14978533333333
$ws.Range("O7").Value = 0.7416121699579786
$ws.Range("P7").Value = 0.7416121699579786
$ws.Range("Q7").Value = 5766.526652439052
$ws.Range("R7").Value = 51898.73987195146
$ws.Range("S7").Value = 0.1164089785754686
$ws.Range("T7").Value = 0.1164089785754686
$ws.Range("I8").Value = 0.1569674599353791
$ws.Range("J8").Value = 0.1569674599353792
$ws.Range("M8").Value = 29.76859933333333
$ws.Range("N8").Value = 89.305798
$ws.Range("O8").Value = 0.2504459365921425
$ws.Range("P8").Value = 0.2504459365921425
$ws.Range("Q8").Value = 1947.383318204557
$ws.Range("R8").Value = 17526.44986384101
$ws.Range("S8").Value = 0.03931186251800563
$ws.Range("T8").Value = 0.03931186251800563
$ws.Range("I9").Value = 0.1569674599353791
$ws.Range("J9").Value = 0.1569674599353792
$ws.Range("M9").Value = 0.2000323333333334
$ws.Range("N9").Value = 0.6000970000000001
$ws.Range("O9").Value = 0.00168289023307462
$ws.Range("P9").Value = 0.00168289023307462
$ws.Range("Q9").Value = 13.08558809479089
$ws.Range("R9").Value = 117.770292853118
$ws.Range("S9").Value = 0.0002641590052357813
$ws.Range("T9").Value = 0.0002641590052357814
$ws.Range("G10").Value = 60.43484133333334
$ws.Range("H10").Value = 181.304524
$ws.Range("I10").Value = 0.1450120099461104
$ws.Range("J10").Value = 0.1450120099461104
$ws.Range("M10").Value = 0.74396
$ws.Range("N10").Value = 2.23188
$ws.Range("O10").Value = 0.006259003216804254
$ws.Range("P10").Value = 0.006259003216804255
$ws.Range("Q10").Value = 44.96110455834667
$ws.Range("R10").Value = 404.64994102512
$ws.Range("S10").Value = 0.0009076306367279552
$ws.Range("T10").Value = 0.0009076306367279555
$ws.Range("G11").Value = 60.43484133333334
$ws.Range("H11").Value = 181.304524
$ws.Range("I11").Value = 0.1450120099461104
$ws.Range("J11").Value = 0.1450120099461104
$ws.Range("M11").Value = 88.14978533333333
$ws.Range("O11").Value = 0.7416121699579786
$ws.Range("P11").Value = 0.7416121699579786
$ws.Range("Q11").Value = 5327.318290187394
$ws.Range("R11").Value = 47945.86461168654
$ws.Range("S11").Value = 0.1075426713661029
$ws.Range("T11").Value = 0.1075426713661029
$ws.Range("G12").Value = 60.43484133333334
$ws.Range("H12").Value = 181.304524
$ws.Range("I12").Value = 0.1450120099461104
$ws.Range("J12").Value = 0.1450120099461104
$ws.Range("M12").Value = 29.76859933333333
$ws.Range("N12").Value = 89.305798
$ws.Range("O12").Value = 0.2504459365921425
$ws.Range("P12").Value = 0.2504459365921425
$ws.Range("Q12").Value = 1799.060577425573
$ws.Range("R12").Value = 16191.54519683015
$ws.Range("S12").Value = 0.03631766864806269
$ws.Range("T12").Value = 0.0363176686480627
$ws.Range("G13").Value = 60.43484133333334
$ws.Range("H13").Value = 181.304524
$ws.Range("I13").Value = 0.1450120099461104
$ws.Range("J13").Value = 0.1450120099461104
$ws.Range("M13").Value = 0.2000323333333334
$ws.Range("N13").Value = 0.6000970000000001
$ws.Range("O13").Value = 0.00168289023307462
$ws.Range("P13").Value = 0.00168289023307462
$ws.Range("Q13").Value = 12.08892232653645
$ws.Range("R13").Value = 108.800300938828
$ws.Range("S13").Value = 0.0002440392952168288
$ws.Range("T13").Value = 0.0002440392952168289
$ws.Range("G14").Value = 89.46554166666668
$ws.Range("H14").Value = 268.396625
$ws.Range("I14").Value = 0.2146705068098712
$ws.Range("J14").Value = 0.2146705068098712
$ws.Range("M14").Value = 0.74396
$ws.Range("N14").Value = 2.23188
$ws.Range("O14").Value = 0.006259003216804254
$ws.Range("P14").Value = 0.006259003216804255
$ws.Range("Q14").Value = 66.55878437833334
$ws.Range("R14").Value = 599.029059405
$ws.Range("S14").Value = 0.001343623392675983
$ws.Range("T14").Value = 0.001343623392675984
$ws.Range("G15").Value = 89.46554166666668
$ws.Range("H15").Value = 268.396625
$ws.Range("I15").Value = 0.2146705068098712
$ws.Range("J15").Value = 0.2146705068098712
$ws.Range("M15").Value = 88.14978533333333
$ws.Range("O15").Value = 0.7416121699579786
$ws.Range("P15").Value = 0.7416121699579786
$ws.Range("Q15").Value = 7886.368292647056
$ws.Range("R15").Value = 70977.31463382349
$ws.Range("S15").Value = 0.1592022603812476
$ws.Range("T15").Value = 0.1592022603812476
$ws.Range("G16").Value = 89.46554166666668
$ws.Range("H16").Value = 268.396625
$ws.Range("I16").Value = 0.2146705068098712
$ws.Range("J16").Value = 0.2146705068098712
$ws.Range("M16").Value = 29.76859933333333
$ws.Range("N16").Value = 89.305798
$ws.Range("O16").Value = 0.2504459365921425
$ws.Range("P16").Value = 0.2504459365921425
$ws.Range("Q16").Value = 2663.263864014639
$ws.Range("R16").Value = 23969.37477613175
$ws.Range("S16").Value = 0.05376335613670809
$ws.Range("T16").Value = 0.0537633561367081
$ws.Range("G17").Value = 89.46554166666668
$ws.Range("H17").Value = 268.396625
$ws.Range("I17").Value = 0.2146705068098712
$ws.Range("J17").Value = 0.2146705068098712
$ws.Range("M17").Value = 0.2000323333333334
$ws.Range("N17").Value = 0.6000970000000001
$ws.Range("O17").Value = 0.00168289023307462
$ws.Range("P17").Value = 0.00168289023307462
$ws.Range("Q17").Value = 17.8960010525139
$ws.Range("R17").Value = 161.064009472625
$ws.Range("S17").Value = 0.000361266899239511
$ws.Range("T17").Value = 0.0003612668992395111
